# Add 4 school questions
# Inserts 4 new "mid class / recess" questions as rows 2-5 (content only; the
# row's Index/Category/Mistakes numbers in columns A/B/G are left untouched,
# matching the source edit which only replaced the Question/Answer text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: bathroom question ---
$ws.Range("C2").Value = "What would you do if you need to go to the bathroom in mid class?"
$ws.Range("D2").Value = "Just get out to the bathroom"
$ws.Range("E2").Value = "Raise my hand and ask the teacher to go to the bathroom"
$ws.Range("F2").Value = "Just shout to the teacher that I'm going to the bathroom"
$ws.Range("H2").Value = "Raise my hand and ask the teacher to go to the bathroom"

# --- Row 3: talking to a friend mid class ---
$ws.Range("C3").Value = "What would you do if you want to talk to your friend mid class?"
$ws.Range("F3").Value = "Talk to him"
$ws.Range("D3").Value = "I'll wait until the end of the class to tell him"
$ws.Range("E3").Value = "I'll pass a note and throw it "
$ws.Range("H3").Value = "I'll wait until the end of the class to tell him"

# --- Row 4: someone threw something at you mid class ---
$ws.Range("C4").Value = "What would you do if someone threw something at you mid class?"
$ws.Range("F4").Value = "I'll tell the teacher to handle it properly"
$ws.Range("E4").Value = "I'll throw at back at them"
$ws.Range("D4").Value = "I will shout ""who threw this at me?!"""
$ws.Range("H4").Value = "I'll tell the teacher to handle it properly"

# --- Row 5: pushed during recess ---
$ws.Range("C5").Value = "In recess someone pushed you, what would you do?"
$ws.Range("D5").Value = "I'll tell him\her it wasn't nice and not to do it again"
$ws.Range("E5").Value = "I'll insult him\her, so he\she won't do it again"
$ws.Range("F5").Value = "I'll fight him\her to make him\her stop"
$ws.Range("H5").Value = "I'll tell him\her it wasn't nice and not to do it again"

# Leave the selection on the last-edited cell, matching the saved workbook view.
[void]$ws.Range("F5").Select()
